$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687460550764"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168749497724"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687494987292"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687495477245"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687496227272"

# Sheet 1 (GNG) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687460161664.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687460373797.csv"
$ws1.Range("B4").Value = "go_stims-1651168746038468.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168746054072.csv"

# Sheet 2 (NB) updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_5-16511687468299563.csv"
$ws2.Range("B3").Value = "OB-16511687468665736.csv"
$ws2.Range("B4").Value = "OB-16511687470775726.csv"
$ws2.Range("B5").Value = "TB-16511687482417235.csv"
$ws2.Range("B6").Value = "ZB-match_7-16511687468001144.csv"
$ws2.Range("B7").Value = "TB-1651168749287758.csv"
$ws2.Range("B8").Value = "ZB-match_8-16511687465175998.csv"
$ws2.Range("B9").Value = "OB-1651168747889727.csv"
$ws2.Range("B10").Value = "TB-16511687494757254.csv"

# Sheet 3 (RS) updates
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687495137262.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168749499727.csv"
$ws4.Range("B4").Value = "MM_stims-16511687495297222.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687495147274.csv"
$ws4.Range("B6").Value = "MM_stims-1651168749545721.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687495307229.csv"

# Sheet 5 (vSAT) updates
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687495517228.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511687495917208.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687496077209.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687495767214.csv"
